# Fix the ordering of fields on the "classFields" sheet for the Order class.
# The rows describing the "private" fields (rows 2-7) need to be reordered;
# row 8 (the "id" field) stays as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# New order (Field Name, Field Type) for rows 2-7, column A (Class Name) and
# column C (Field Modifier) are unchanged for all of these rows.
$newFields = @(
    @("status",       "java.lang.String"),
    @("productCount",  "int"),
    @("source",        "java.lang.String"),
    @("customerId",    "java.lang.Long"),
    @("productId",     "java.lang.Long"),
    @("price",         "int")
)

for ($i = 0; $i -lt $newFields.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $newFields[$i][0]
    $ws.Cells.Item($row, 4).Value = $newFields[$i][1]
}
